$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: highlight the "Today,Tomorrow" note cells (D12:F12) ---
# D12 and E12 keep their existing text but get bold-red-on-yellow formatting
# with the existing full border; F12 gets a brand new value with the same
# bold-red-on-yellow formatting but no border.
$ws.Range("D12").Font.Bold = $true
$ws.Range("D12").Font.Color = 255
$ws.Range("D12").Interior.Color = 65535

$ws.Range("E12").Font.Bold = $true
$ws.Range("E12").Font.Color = 255
$ws.Range("E12").Interior.Color = 65535

$ws.Range("F12").Value = "Today,Tomorrow"
$ws.Range("F12").Font.Bold = $true
$ws.Range("F12").Font.Color = 255
$ws.Range("F12").Interior.Color = 65535

# --- Row 7: fill in the new names ---
# A7 should pick up the same vertical-rule-only style already used by E2.
$ws.Range("E2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = "Ram "
$ws.Range("B7").Value = "Maneesh"
$ws.Range("F7").Value = "goutham Allu"

# --- selection cursor moves to F7 ---
$null = $ws.Range("F7").Select()
